# Elimna EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Periodo Mora" column (E16:E21) lists the overdue payroll periods.
# The database is refreshed: the newest period (1701) now appears first
# and the oldest (1608) last -- i.e. the six period codes are written in
# reverse chronological order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "1701"
$ws.Range("E17").Value = "1612"
$ws.Range("E18").Value = "1611"
$ws.Range("E19").Value = "1610"
$ws.Range("E20").Value = "1609"
$ws.Range("E21").Value = "1608"
